$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.87'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.51%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '9.36%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.166'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.27%'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.46%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.589'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.25%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8566'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '4.59%'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8669'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.16%'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1365'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2.63%'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07065'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.71%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02931'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '3.45%'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09380'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.14%'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001516'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-1.01%'
$ws.Range('B14').Value = 'CoinExToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04135'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '2.25%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0005997'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-94.09%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006110'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.65%'
$ws.Range('B17').Value = 'UpBots'
$ws.Range('C17').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.007489'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '5,070.96%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.489'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.61%'
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.098'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.97%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.67%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.56%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.03399'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '5.51%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.1300'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.10%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.469'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-3.42%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.005009'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '11.98%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '0.53%'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '22.19%'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.55%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005726'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-4.37%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1070'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '1.15%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002426'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.48%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008515'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-12.34%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005251'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.70%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.03%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '1.10%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.03%'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.03%'
